# This document has several paragraphs whose runs begin with an empty
# <w:r/> run immediately followed by the text-bearing run. A plain
# Find.Execute(...,Replace:=2) ends up merging/dropping that leading
# empty run because the edit touches position 0 of the paragraph.
# To faithfully reproduce the target diff (which keeps those empty runs
# untouched) we replace paragraph text surgically with Range.InsertXML,
# targeting only the run that holds the text (start of paragraph through
# the char just before the paragraph mark) and rebuilding just that run
# (with its original direct formatting, if any). That leaves any sibling
# runs - such as the leading empty run - completely undisturbed.

function Replace-ParagraphText {
    param(
        [__ComObject]$doc,
        [int]$paraIndex,
        [string]$newText,
        [string]$rPrXml = ""
    )

    $p = $doc.Paragraphs($paraIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    # Range covering the paragraph's text but excluding the trailing
    # paragraph-mark character.
    $rng = $doc.Range($pStart, $pEnd - 1)

    $escaped = $newText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    if ($newText -ne $newText.Trim()) {
        $tOpen = '<w:t xml:space="preserve">'
    } else {
        $tOpen = '<w:t>'
    }

    $runXml = '<w:r>' + $rPrXml + $tOpen + $escaped + '</w:t></w:r>'

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
}

$d = $word.ActiveDocument

# 1. Title heading (Heading1 at the very top of the document).
Replace-ParagraphText $d 1 "Play Chunky Fruits Online for Free"

# 2. "What we like" bullet list rewording/reordering.
Replace-ParagraphText $d 42 "Modern gameplay experience with a 5x3 game grid"
Replace-ParagraphText $d 43 "Expanding Wild symbols trigger free spins"
Replace-ParagraphText $d 44 "Immersive and nostalgic atmosphere"

# 3. "What we don't like" bullet list rewording.
Replace-ParagraphText $d 47 "Limited variety of symbols"
Replace-ParagraphText $d 48 "No progressive jackpot feature"

# 4. Bold title repeated near the end of the document.
Replace-ParagraphText $d 49 "Play Chunky Fruits Online for Free" "<w:rPr><w:b/></w:rPr>"

# 5. Italic meta description near the end of the document.
Replace-ParagraphText $d 50 "Read our review of Chunky Fruits and play this slot game for free. Enjoy modern gameplay and expanding Wild symbols." "<w:rPr><w:i/></w:rPr>"
